$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2: Area column starts here, plus Atotal/H2, and summary columns J2/K2
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3: first explicit (non-shared) Area formula
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15: shared Area formula pattern (D[r]-D[r-1])*B[r]/100
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 7).Formula = "=(D$r-D$prev)*B$r/100"
}

# Match the author's final selection on the new summary cells
$ws.Range("J2:K2").Select()
